$d = $word.ActiveDocument

# 1) Merge "3 - Pending" + " " into a single run's text "3 - Pending "
#    (visible text is unchanged: "3 - Pending" followed by a space)
$d.Content.Find.Execute("3 - Pending ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3 - Pending ", 2) | Out-Null

# 2) Change "Bachelor of Science in Geography" to "Bachelor of Art in Geography"
$d.Content.Find.Execute("Bachelor of Science in Geography", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Bachelor of Art in Geography", 2) | Out-Null
